$wb = $excel.ActiveWorkbook

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 132
$ws.Range("H132").Value = 1810.25
$ws.Range("I132").Value = 1446
$ws.Range("J132").Value = 2356.625
$ws.Range("K132").Value = 4338
$ws.Range("L132").Value = 7069.875
$ws.Range("M132").Value = -1808
$ws.Range("N132").Value = -12129.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 63
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 66
$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864

# Row 75
$ws.Range("H75").Value = 11500
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 18000
$ws.Range("K75").Value = 5000
$ws.Range("L75").Value = 18000
$ws.Range("M75").Value = -4064
$ws.Range("N75").Value = -19872

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 78
$ws.Range("H78").Value = 11500
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 18000
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 54000
$ws.Range("M78").Value = -10320
$ws.Range("N78").Value = -63360

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 82
$ws.Range("H82").Value = 16955.75
$ws.Range("I82").Value = 6628.5
$ws.Range("J82").Value = 27283
$ws.Range("K82").Value = 6628.5
$ws.Range("L82").Value = 27283
$ws.Range("M82").Value = -6245.5
$ws.Range("N82").Value = -28049

# Row 85
$ws.Range("H85").Value = 16955.75
$ws.Range("I85").Value = 6628.5
$ws.Range("J85").Value = 27283
$ws.Range("K85").Value = 6628.5
$ws.Range("L85").Value = 27283
$ws.Range("M85").Value = -5302.5
$ws.Range("N85").Value = -29935

# Row 86
$ws.Range("H86").Value = 35719024
$ws.Range("I86").Value = 76925600
$ws.Range("J86").Value = 6660.4
$ws.Range("K86").Value = 76925600
$ws.Range("L86").Value = 6660.4
$ws.Range("M86").Value = -76924477
$ws.Range("N86").Value = -8906.4

# Row 88
$ws.Range("H88").Value = 19000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 19000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 19000
$ws.Range("N88").Value = -19812

# Row 89
$ws.Range("H89").Value = 35719024
$ws.Range("I89").Value = 76925600
$ws.Range("J89").Value = 6660.4
$ws.Range("K89").Value = 384628000
$ws.Range("L89").Value = 33302
$ws.Range("M89").Value = -384622384
$ws.Range("N89").Value = -44534

# Row 91
$ws.Range("H91").Value = 19000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 19000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 19000
$ws.Range("N91").Value = -21808

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 96
$ws.Range("H96").Value = 11594.066
$ws.Range("I96").Value = 7909.25
$ws.Range("J96").Value = 26333.334
$ws.Range("K96").Value = 7909.25
$ws.Range("L96").Value = 26333.334
$ws.Range("M96").Value = -5163.25
$ws.Range("N96").Value = -31825.334

# Row 97
$ws.Range("H97").Value = 4463.778
$ws.Range("I97").Value = 4084.25
$ws.Range("J97").Value = 7500
$ws.Range("K97").Value = 4084.25
$ws.Range("L97").Value = 7500
$ws.Range("M97").Value = -3093.25
$ws.Range("N97").Value = -9482

# Row 98
$ws.Range("H98").Value = 30000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 30000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

# Row 99
$ws.Range("H99").Value = 2070
$ws.Range("I99").Value = 1322.25
$ws.Range("J99").Value = 2568.5
$ws.Range("K99").Value = 1322.25
$ws.Range("L99").Value = 2568.5
$ws.Range("M99").Value = 175.75
$ws.Range("N99").Value = -5564.5

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 102
$ws.Range("H102").Value = 12764
$ws.Range("I102").Value = 12764
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 12764
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9519

# Row 103
$ws.Range("H103").Value = 23324.545
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 23324.545
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 23324.545
$ws.Range("N103").Value = -25668.545

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105
$ws.Range("H105").Value = 1696.6296
$ws.Range("I105").Value = 1696.2084
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1696.2084
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 50.79160000000002
$ws.Range("N105").Value = -5194

# Row 106
$ws.Range("H106").Value = 26330
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 26330
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 26330
$ws.Range("N106").Value = -28854

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 112
$ws.Range("H112").Value = 48000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 48000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 48000
$ws.Range("N112").Value = -50954

# Row 114
$ws.Range("H114").Value = 40674
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 40674
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 40674
$ws.Range("N114").Value = -49352

# Row 117
$ws.Range("H117").Value = 40000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 40000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

# Row 121
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
